$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ "B"=1.02; "C"=1.031825702612842; "D"=1.041311234128588; "E"=1.031347056317505; "F"=1.0498088441016; "I"=1.037609870213082; "J"=1.03695929352739; "K"=1.044091108366749; "L"=1.034155437614137; "M"=1.05256487077546 }
    3 = @{ "B"=1.02; "C"=1.032826317500215; "D"=1.042099433132061; "E"=1.032198089628791; "F"=1.050730073968098; "I"=1.037832511416765; "J"=1.03760168606342; "K"=1.044689967995596; "L"=1.03481490173521; "M"=1.0532981399888 }
    4 = @{ "B"=1.02; "C"=1.033474039065305; "D"=1.042609314512821; "E"=1.032749369543917; "F"=1.051326352033497; "I"=1.037974874079927; "J"=1.03801702990991; "K"=1.045076679039764; "L"=1.035241582643736; "M"=1.053772165986827 }
    5 = @{ "B"=1.02; "C"=1.033746402167627; "D"=1.042823634939292; "E"=1.032981271293592; "F"=1.051577069487032; "I"=1.038034315530487; "J"=1.03819156142429; "K"=1.045239062078493; "L"=1.035420949772614; "M"=1.053971338238057 }
    6 = @{ "B"=1.02; "C"=1.033792136690164; "D"=1.042859618272524; "E"=1.033020217012734; "F"=1.051619168490187; "I"=1.038044272080601; "J"=1.038220861377254; "K"=1.045266315721368; "L"=1.035451065721163; "M"=1.054004773757534 }
    7 = @{ "B"=1.02; "C"=1.033477678155551; "D"=1.042612178406718; "E"=1.032752467663212; "F"=1.051329701968152; "I"=1.037975669942647; "J"=1.038019362320044; "K"=1.04507884955867; "L"=1.035243979394519; "M"=1.05377482776289 }
    8 = @{ "B"=1.02; "C"=1.032163811604543; "D"=1.041577637171714; "E"=1.0316345412319; "F"=1.050120139955654; "I"=1.037685464728975; "J"=1.037176460375553; "K"=1.044293659077148; "L"=1.034378313677912; "M"=1.052812774760596 }
    9 = @{ "B"=1.02; "C"=1.029850600772034; "D"=1.039753650527361; "E"=1.029669286830516; "F"=1.047990170215743; "I"=1.037161082700745; "J"=1.035688686554784; "K"=1.042904031839744; "L"=1.032852655458058; "M"=1.051114124289294 }
    10 = @{ "B"=1.02; "C"=1.028309826647406; "D"=1.038537056918057; "E"=1.02836232008268; "F"=1.046571217865053; "I"=1.036802784554936; "J"=1.034695216552274; "K"=1.041973615309906; "L"=1.031835430989561; "M"=1.049979467088316 }
    11 = @{ "B"=1.02; "C"=1.027642983782997; "D"=1.03801012874233; "E"=1.027797160882846; "F"=1.045957052438766; "I"=1.036645578124839; "J"=1.034264656989911; "K"=1.041569796377041; "L"=1.031394942503397; "M"=1.04948763182629 }
    12 = @{ "B"=1.02; "C"=1.027395337433533; "D"=1.037814384860605; "E"=1.027587351611615; "F"=1.045728962632532; "I"=1.036586875492513; "J"=1.034104671444994; "K"=1.041419659201009; "L"=1.031231322611234; "M"=1.049304864843533 }
    13 = @{ "B"=1.02; "C"=1.027448456203374; "D"=1.037856373444813; "E"=1.027632351147447; "F"=1.045777886908879; "I"=1.036599481389203; "J"=1.034138991436784; "K"=1.041451870496483; "L"=1.031266419751742; "M"=1.049344072473904 }
    14 = @{ "B"=1.02; "C"=1.027622512261454; "D"=1.037993948879844; "E"=1.027779815610894; "F"=1.045938197671618; "I"=1.036640732051077; "J"=1.034251433684985; "K"=1.041557388860114; "L"=1.031381417686358; "M"=1.049472525825961 }
    15 = @{ "B"=1.02; "C"=1.027729760417333; "D"=1.038078711118991; "E"=1.027870688735297; "F"=1.046036975564579; "I"=1.036666106992092; "J"=1.034320705589655; "K"=1.041622383565505; "L"=1.031452271351395; "M"=1.049551659935468 }
    16 = @{ "B"=1.02; "C"=1.028354089634922; "D"=1.038572024644756; "E"=1.028399844056421; "F"=1.046611983358499; "I"=1.036813174446118; "J"=1.034723783408128; "K"=1.042000395650258; "L"=1.031864664321567; "M"=1.050012097665223 }
    17 = @{ "B"=1.02; "C"=1.028745801433265; "D"=1.038881431722434; "E"=1.028731974955534; "F"=1.046972738258356; "I"=1.036904874529284; "J"=1.034976522038083; "K"=1.042237260855183; "L"=1.032123341709371; "M"=1.050300779139543 }
    18 = @{ "B"=1.02; "C"=1.028974311392333; "D"=1.039061890532466; "E"=1.028925775137868; "F"=1.047183184573739; "I"=1.036958162751221; "J"=1.035123903451707; "K"=1.042375329340589; "L"=1.032274221532272; "M"=1.050469111814577 }
    19 = @{ "B"=1.02; "C"=1.02905223264727; "D"=1.039123420070001; "E"=1.028991868489305; "F"=1.047254945384561; "I"=1.036976298909023; "J"=1.03517415045243; "K"=1.042422391671365; "L"=1.032325667238201; "M"=1.050526500341752 }
    20 = @{ "B"=1.02; "C"=1.028703771245716; "D"=1.038848236619251; "E"=1.028696332808068; "F"=1.046934030173339; "I"=1.036895056546368; "J"=1.034949409372057; "K"=1.042211856865522; "L"=1.032095588308807; "M"=1.050269811550681 }
    21 = @{ "B"=1.02; "C"=1.027571255741319; "D"=1.037953436867826; "E"=1.027736387791688; "F"=1.045890989074871; "I"=1.036628593289973; "J"=1.034218323801665; "K"=1.041526320214021; "L"=1.031347553734187; "M"=1.049434701647166 }
    22 = @{ "B"=1.02; "C"=1.026859480693635; "D"=1.037390729603625; "E"=1.027133503971444; "F"=1.045235411422278; "I"=1.036459268766153; "J"=1.033758334391626; "K"=1.041094481597218; "L"=1.030877218524364; "M"=1.048909187114682 }
    23 = @{ "B"=1.02; "C"=1.027236779235031; "D"=1.037689041572637; "E"=1.027453040110325; "F"=1.045582924013206; "I"=1.036549200271796; "J"=1.034002214263132; "K"=1.04132348437715; "L"=1.031126553423345; "M"=1.049187814478003 }
    24 = @{ "B"=1.02; "C"=1.028722762782149; "D"=1.038863236098051; "E"=1.028712437732181; "F"=1.046951520616402; "I"=1.036899493484767; "J"=1.034961660532346; "K"=1.042223336116082; "L"=1.032108128883081; "M"=1.050283804626682 }
    25 = @{ "B"=1.02; "C"=1.030448382467065; "D"=1.040225305730963; "E"=1.03017679204069; "F"=1.048540642014826; "I"=1.037298185619009; "J"=1.036073600478375; "K"=1.04326399217456; "L"=1.033247098970678; "M"=1.051553662137748 }
}

foreach ($rowKey in $data.Keys) {
    $rowNum = [int]$rowKey
    $rowValues = $data[$rowKey]
    foreach ($col in $rowValues.Keys) {
        $ws.Range("$col$rowNum").Value = $rowValues[$col]
    }
}
